$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. Make the rubric table span the full page width (100%).
$t.PreferredWidthType = 2   # wdPreferredWidthPercent
$t.PreferredWidth = 250     # -> <w:tblW w:type="pct" w:w="5000"/>

# 2. Fill in the three empty rubric cells in the data row (row 2) with
#    left-justified text describing each performance level.
$cellXmlTemplate = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">{0}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$texts = @(
    "Response directly addresses the prompt with specific details from the readings and NotebookLM. Includes concrete examples from field experience or teaching practice.",
    "Response addresses the prompt but lacks specific details or examples. May be vague or general.",
    "No response or response does not address the prompt."
)

for ($i = 1; $i -le 3; $i++) {
    $cell = $t.Cell(2, $i)
    $xml = $cellXmlTemplate -f $texts[$i - 1]
    [void]$cell.Range.InsertXML($xml)
}
